$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared text labels ---
$ws.Range("A4").Value = "Stream-sugarcane"
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"

# --- Clear column D data (D4:D13) ---
$ws.Range("D4:D13").ClearContents()

# --- Update numeric values for rows 4-13 ---
$ws.Range("C4").Value = 0.03827795694892373
$ws.Range("E4").Value = -0.06061201530038251
$ws.Range("F4").Value = -0.06181354533863348
$ws.Range("H4").Value = -0.06061201530038251
$ws.Range("I4").Value = -0.08235205880147006
$ws.Range("J4").Value = -0.1558450909518538

$ws.Range("C5").Value = 0.02411910297757444
$ws.Range("E5").Value = 0.1472931823295582
$ws.Range("F5").Value = 0.1326213155328883
$ws.Range("H5").Value = 0.1412060301507538
$ws.Range("I5").Value = 0.08803720093002328
$ws.Range("J5").Value = 0.01807414841114307

$ws.Range("C6").Value = -0.003130578264456612
$ws.Range("E6").Value = -0.05566489162229056
$ws.Range("F6").Value = -0.06832970824270607
$ws.Range("H6").Value = -0.05884647116177905
$ws.Range("I6").Value = -0.06873471836795922
$ws.Range("J6").Value = -0.004111647505597405

$ws.Range("C7").Value = 0.083625590639766
$ws.Range("E7").Value = 0.7678016950423762
$ws.Range("F7").Value = 0.7849681242031051
$ws.Range("H7").Value = 0.7736998424960625
$ws.Range("I7").Value = 0.9997449936248407
$ws.Range("J7").Value = 0.2469411090268892

$ws.Range("C8").Value = 0.9619410485262132
$ws.Range("E8").Value = 0.03280582014550364
$ws.Range("F8").Value = 0.02794869871746794
$ws.Range("H8").Value = 0.03044326108152704
$ws.Range("I8").Value = 0.04204005100127504
$ws.Range("J8").Value = 0.02033172940199461

$ws.Range("C9").Value = -0.01385134628365709
$ws.Range("E9").Value = -0.02015900397509938
$ws.Range("F9").Value = -0.02118802970074252
$ws.Range("H9").Value = -0.01914647866196655
$ws.Range("I9").Value = 0.02138603465086627
$ws.Range("J9").Value = -0.007806280050758446

$ws.Range("C10").Value = 0.113422335558389
$ws.Range("E10").Value = -0.002806570164254107
$ws.Range("F10").Value = -0.01432985824645617
$ws.Range("H10").Value = -0.00570014250356259
$ws.Range("I10").Value = -0.03298732468311708
$ws.Range("J10").Value = 0.04264652994678374

$ws.Range("C11").Value = 0.04565964149103728
$ws.Range("E11").Value = 0.1034230855771394
$ws.Range("F11").Value = 0.09427285682142056
$ws.Range("H11").Value = 0.102137553438836
$ws.Range("I11").Value = 0.1267321683042076
$ws.Range("J11").Value = 0.04857549264867127

$ws.Range("C12").Value = 0.0834455861396535
$ws.Range("E12").Value = 0.5355283882097053
$ws.Range("F12").Value = 0.5153573839345984
$ws.Range("H12").Value = 0.5315637890947275
$ws.Range("I12").Value = -0.09208580214505364
$ws.Range("J12").Value = -0.1507861594623875

$ws.Range("C13").Value = -0.2088097202430061
$ws.Range("E13").Value = 0.0003735093377334434
$ws.Range("F13").Value = 0.007227180679516989
$ws.Range("H13").Value = 0.001819545488637216
$ws.Range("I13").Value = 0.06154353858846472
$ws.Range("J13").Value = -0.04711894039477304

